# Rename "Examples Applications" to "Example Applications" on the single
# slide of the networking_overview deck.
#
# The target text box (shape "Rectangle 48") holds the caption as one run:
#   "Examples Applications"
# and needs to become two runs:
#   "Example " + "Applications"
# (i.e. just the trailing "s" of "Examples" is removed), matching the
# author's commit "rename Examples Applications to Example Applications".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the shape whose text is exactly "Examples Applications" (robust to
# shape re-ordering / re-numbering).
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "Examples Applications") {
            $targetShape = $shp
            break
        }
    }
}

if ($targetShape -eq $null) {
    Write-Host "Could not find target shape with text 'Examples Applications'"
} else {
    $tr = $targetShape.TextFrame.TextRange

    # Insert "Example " at the very start; this becomes its own run ahead
    # of the existing "Examples Applications" run.
    $null = $tr.Characters(1, 0).InsertBefore("Example ")

    # The text is now "Example Examples Applications" - remove the old
    # leading "Examples " (9 characters) that immediately follows what we
    # just inserted, leaving "Example Applications" split as
    # "Example " + "Applications".
    $tr2 = $targetShape.TextFrame.TextRange
    $tr2.Characters(9, 9).Text = ""

    Write-Host "Updated text: $($targetShape.TextFrame.TextRange.Text)"
}
